$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.892.40'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.901.97'
$ws.Range('E3').Value = '  +3.02%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '351.65'
$ws.Range('E5').Value = '  -0.10%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.81'
$ws.Range('E6').Value = '  -0.33%  '

$ws.Range('E7').Value = '  -0.80%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.616'
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.23'
$ws.Range('E10').Value = '  -2.94%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0873'
$ws.Range('E11').Value = '  +3.71%  '

$ws.Range('E12').Value = '  +0.67%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.92'
$ws.Range('E13').Value = '  +0.01%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.363.78'
$ws.Range('E14').Value = '  +3.36%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.68'
$ws.Range('E15').Value = '  -1.77%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.908.47'
$ws.Range('E16').Value = '  +2.98%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.977'
$ws.Range('E17').Value = '  +0.98%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.890.93'
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('E19').Value = '  -1.51%  '

$ws.Range('E20').Value = '  -2.76%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.03'
$ws.Range('E21').Value = '  +3.08%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  -0.48%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.81'
$ws.Range('E23').Value = '  +0.37%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.63'
$ws.Range('E24').Value = '  -0.67%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.77'
$ws.Range('E25').Value = '  +0.11%  '

$ws.Range('E26').Value = '  +9.02%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '26.63'
$ws.Range('E27').Value = '  +1.47%  '

$ws.Range('E28').Value = '  +0.06%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.85'
$ws.Range('E29').Value = '  +10.93%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '10.55'
$ws.Range('E30').Value = '  +0.47%  '

$ws.Range('E31').Value = '  +12.56%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '36.80'
$ws.Range('E32').Value = '  -5.65%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.96'
$ws.Range('E33').Value = '  +4.50%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '52.87'
$ws.Range('E34').Value = '  +0.16%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.09'
$ws.Range('E35').Value = '  -7.98%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0450'
$ws.Range('E36').Value = '  -0.70%  '

$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  +3.61%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.52'
$ws.Range('E39').Value = '  -2.69%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  +0.15%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.67'
$ws.Range('E41').Value = '  +5.05%  '

$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '22.73'
$ws.Range('E43').Value = '  +2.15%  '

$ws.Range('E44').Value = '  -2.18%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.52'
$ws.Range('E45').Value = '  +2.10%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.185.38'
$ws.Range('E46').Value = '  +2.48%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.47'
$ws.Range('E47').Value = '  -1.84%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '110.28'
$ws.Range('E48').Value = '  -9.53%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.249'
$ws.Range('E49').Value = '  +10.44%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0340'
$ws.Range('E50').Value = '  +4.94%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.937'
$ws.Range('E51').Value = '  -7.44%  '
